$p = $ppt.ActivePresentation
try {
  $d2 = $p.Designs.Item(2)
  Write-Output ("Design2: " + $d2.Name)
} catch {
  Write-Output ("ERR: " + $_.Exception.Message)
}
